$p = $ppt.ActivePresentation

# Slide 4 - "Title 1" holds the plain (unstyled) sprint-goals summary list,
# one run per line separated by manual line breaks. Edit back-to-front so
# that earlier (still-unedited) runs keep their original character offsets.
$s4 = $p.Slides.Item(4)
$titleRange = $s4.Shapes.Item(1).TextFrame.TextRange

$titleRange.Characters(100, 21).Text = "- update user experience"
$titleRange.Characters(43, 26).Text = "- Route schedule On Home page"
$titleRange.Characters(17, 25).Text = "- Emails getting Sent"
$titleRange.Characters(1, 15).Text = "- Added Sign in directory"

# Slide 6 - "Text Placeholder 2" holds the same list but styled (red,
# all-caps). Only the "Route Modal" line changed here.
$s6 = $p.Slides.Item(6)
$bodyRange = $s6.Shapes.Item(2).TextFrame.TextRange

$bodyRange.Characters(43, 26).Text = "- Route schedule On Home page"
